$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Courses (Registration, Deletion)." -> "Courses (Registration, Deletion/Archive)."
#    Insert "/Archive" right before the closing paren, leaving the surrounding
#    "Courses" run and the trailing "." run untouched.
# ---------------------------------------------------------------------------
$marker = " (Registration, Deletion)"
$full = $d.Content.Text
$idx = $full.IndexOf($marker)
if ($idx -ge 0) {
    $insertAt = $idx + $marker.Length - 1   # position right before the ")"
    $gap = $d.Range($insertAt, $insertAt)
    $gap.InsertBefore("/Archive")
}

# ---------------------------------------------------------------------------
# 2. Insert a new "Course Read-only." bullet paragraph right after
#    "Course Student/Instructor Registration." (i.e. right before "General Feeds.")
# ---------------------------------------------------------------------------
foreach ($para in @($d.Paragraphs)) {
    $t = $para.Range.Text.Trim()
    if ($t -eq "Course Student/Instructor Registration.") {
        $para.Range.InsertParagraphAfter()
        $inserted = $para.Next()
        $body = $inserted.Range
        $body.Collapse(1)
        $body.InsertAfter("Course Read-only")
        $body.Collapse(0)
        $body.InsertAfter(".")
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Remove the old "Course Read-only and Archive." bullet paragraph
#    (it used to sit right after "Send a confirmation mail after submitting
#    an assignment.", just before "Support Files (file, audio, video, and link).")
# ---------------------------------------------------------------------------
foreach ($para in @($d.Paragraphs)) {
    $t = $para.Range.Text.Trim()
    if ($t -eq "Course Read-only and Archive.") {
        $para.Range.Delete()
        break
    }
}
